# Atualização automática de preços de eletricidade
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = 46041
$ws.Range("B2").Value = 119.35
$ws.Range("C2").Value = 110.01
$ws.Range("D2").Value = 102.21
$ws.Range("E2").Value = 100.44
$ws.Range("F2").Value = 100.14
$ws.Range("G2").Value = 102.94
$ws.Range("H2").Value = 126.97
$ws.Range("I2").Value = 139.96
$ws.Range("J2").Value = 173.69
$ws.Range("K2").Value = 151.07
$ws.Range("L2").Value = 132.29
$ws.Range("M2").Value = 121.91
$ws.Range("N2").Value = 119.24
$ws.Range("O2").Value = 118.99
$ws.Range("P2").Value = 111.96
$ws.Range("Q2").Value = 115.45
$ws.Range("R2").Value = 123.08
$ws.Range("S2").Value = 138.92
$ws.Range("T2").Value = 183.15
$ws.Range("U2").Value = 181.57
$ws.Range("V2").Value = 149.62
$ws.Range("W2").Value = 135.34
$ws.Range("X2").Value = 126.58
$ws.Range("Y2").Value = 112.62
$ws.Range("Z2").Value = 129.06
$ws.Range("AA2").Value = "16h-20h"
$ws.Range("AB2").Value = 156.68
$ws.Range("AC2").Value = "18h-20h"
$ws.Range("AD2").Value = 182.36
$ws.Range("AE2").Value = "8h-10h"
$ws.Range("AF2").Value = 162.38
$ws.Range("AG2").Value = "0h-23h"
